$wb = $excel.ActiveWorkbook

# Rename sheets with updated timestamp-based identifiers
$wb.Worksheets.Item(1).Name = "GNG_TO-1650291198861876"
$wb.Worksheets.Item(2).Name = "NB_TO-16502912041034446"
$wb.Worksheets.Item(3).Name = "RS_TO-16502912041053865"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502912041664267"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1650291204243387"

# Sheet 1 (GNG) - update stim file names
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502911988338752.csv"
$ws1.Range("B3").Value = "GNG_stims-1650291198844903.csv"
$ws1.Range("B4").Value = "go_stims-16502911988468711.csv"
$ws1.Range("B5").Value = "GNG_stims-1650291198860906.csv"

# Sheet 2 (NB) - update stim file names
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_7-16502911990428762.csv"
$ws2.Range("B3").Value = "OB-16502911994673865.csv"
$ws2.Range("B4").Value = "TB-1650291204086387.csv"
$ws2.Range("B5").Value = "ZB-match_6-16502911991713846.csv"
$ws2.Range("B6").Value = "OB-16502912002433894.csv"
$ws2.Range("B7").Value = "OB-16502911998003867.csv"
$ws2.Range("B8").Value = "TB-16502912007123876.csv"
$ws2.Range("B9").Value = "TB-16502912030863855.csv"
$ws2.Range("B10").Value = "ZB-match_4-16502911988679056.csv"

# Sheet 4 (TOL) - update stim file names
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-1650291204119387.csv"
$ws4.Range("B3").Value = "ZM_stims-16502912041073887.csv"
$ws4.Range("B4").Value = "MM_stims-16502912041503859.csv"
$ws4.Range("B5").Value = "ZM_stims-16502912041203868.csv"
$ws4.Range("B6").Value = "MM_stims-1650291204165389.csv"
$ws4.Range("B7").Value = "ZM_stims-16502912041513872.csv"

# Sheet 5 (vSAT) - update stim file names
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16502912041703854.csv"
$ws5.Range("B3").Value = "vSAT_stims-16502912042114182.csv"
$ws5.Range("B4").Value = "SAT_stims-1650291204196385.csv"
$ws5.Range("B5").Value = "vSAT_stims-16502912042273886.csv"
